$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{Row=5; I='sd'; J='Statement-non-opinion'}
    @{Row=8; I='sd'; J='Statement-non-opinion'}
    @{Row=14; I='sv'; J='Statement-opinion'}
    @{Row=15; I='sd'; J='Statement-non-opinion'}
    @{Row=16; I='sv'; J='Statement-opinion'}
    @{Row=21; I='sd'; J='Statement-non-opinion'}
    @{Row=24; I='sv'; J='Statement-opinion'}
    @{Row=25; I='aa'; J='Agree/Accept'}
    @{Row=42; I='%'; J='Uninterpretable'}
    @{Row=90; I='aa'; J='Agree/Accept'}
    @{Row=94; I='qy'; J='Yes-No-Question'}
    @{Row=113; I='sv'; J='Statement-opinion'}
    @{Row=115; I='sv'; J='Statement-opinion'}
    @{Row=121; I='aa'; J='Agree/Accept'}
    @{Row=127; I='sv'; J='Statement-opinion'}
    @{Row=128; I='sv'; J='Statement-opinion'}
    @{Row=146; I='ba'; J='Appreciation'}
    @{Row=158; I='aa'; J='Agree/Accept'}
    @{Row=167; I='aa'; J='Agree/Accept'}
    @{Row=175; I='aa'; J='Agree/Accept'}
    @{Row=179; I='aa'; J='Agree/Accept'}
    @{Row=191; I='ba'; J='Appreciation'}
    @{Row=194; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=199; I='aa'; J='Agree/Accept'}
    @{Row=201; I='ba'; J='Appreciation'}
    @{Row=205; I='aa'; J='Agree/Accept'}
    @{Row=232; I='sd'; J='Statement-non-opinion'}
    @{Row=235; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=246; I='sd'; J='Statement-non-opinion'}
    @{Row=263; I='sd'; J='Statement-non-opinion'}
    @{Row=269; I='sv'; J='Statement-opinion'}
    @{Row=276; I='sd'; J='Statement-non-opinion'}
    @{Row=283; I='sv'; J='Statement-opinion'}
    @{Row=295; I='sd'; J='Statement-non-opinion'}
    @{Row=300; I='aa'; J='Agree/Accept'}
    @{Row=305; I='sd'; J='Statement-non-opinion'}
    @{Row=316; I='sv'; J='Statement-opinion'}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}

$wb.Save()